# feat: add stock search functionality and enhance stock management UI
#
# The stock list was refreshed: the "creation_date" timestamp (column G) for
# the pre-existing rows gets normalised to a single, shared refresh
# timestamp, and a brand-new product row (row 6 - "Chitos paquete grande")
# is appended to the bottom of the table, as it would be when a user adds a
# new stock item through the UI.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Normalise the "creation_date" column for the existing rows -----------
# All four existing products now carry the same refreshed timestamp.
$refreshedTimestamp = 45803.7748408449

$ws.Range("G2").Value = $refreshedTimestamp
$ws.Range("G3").Value = $refreshedTimestamp
$ws.Range("G4").Value = $refreshedTimestamp
$ws.Range("G5").Value = $refreshedTimestamp

# --- Append the new product row --------------------------------------------
$newRow = 6

$ws.Cells.Item($newRow, 1).Value = "5AYB"
$ws.Cells.Item($newRow, 2).Value = "Alimentos y bebidas"
$ws.Cells.Item($newRow, 3).Value = "Chitos paquete grande"

# quantity / purchase_price / sale_price arrive as plain (string) form input
# from the "add stock" UI, so they are stored as text, same as the other
# textual columns.
$ws.Range("D6:F6").NumberFormat = "@"
$ws.Cells.Item($newRow, 4).Value = "50"
$ws.Cells.Item($newRow, 5).Value = "10000"
$ws.Cells.Item($newRow, 6).Value = "12000"

# creation_date for the new row is "now" (later than the refreshed
# timestamp above), formatted the same way as the rest of column G.
$ws.Range("G6").NumberFormat = $ws.Range("G2").NumberFormat()
$ws.Cells.Item($newRow, 7).Value = 45803.80732169651
